$d = $word.ActiveDocument

# The document has three tables: UC1 "Registrere bruker" (kept), UC2 "Bestille
# Broyting" (removed) and UC3 "Motta og behandle broyteforesporsel" (removed).
# Capture every offset we need up front (before any deletion shifts them),
# then apply the deletions from the end of the document backwards so earlier
# offsets stay valid.

$strayEnd   = $d.Tables(2).Range.Start   # right after the lone stray empty paragraph
$gapStart   = $d.Tables(2).Range.End     # spacer paragraphs between UC2 and UC3
$tailStart  = $d.Tables(3).Range.End     # heading + 2 empty paragraphs after UC3

# Drop the bold heading paragraph plus two empty paragraphs that sat right
# after the UC3 table, leaving the rest of the trailing empty paragraphs
# untouched.
$d.Range($tailStart, $tailStart + 3).Delete()

# Remove the whole UC3 table.
$d.Tables(3).Delete()

# Drop the spacer paragraphs (seven empty + one bold heading paragraph) that
# separated the UC2 and UC3 tables.
$d.Range($gapStart, $gapStart + 8).Delete()

# Remove the whole UC2 table.
$d.Tables(2).Delete()

# Drop the single stray empty paragraph that sat between the "   " paragraph
# (right after the UC1 table) and the start of the (now removed) UC2 table.
$d.Range($strayEnd - 1, $strayEnd).Delete()
